$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NumberError")

$values = @(
  "Tue Nov 12 17:42:17 EST 2024",
  "Tue Nov 12 17:42:30 EST 2024",
  "Tue Nov 12 17:42:39 EST 2024",
  "Tue Nov 12 17:42:48 EST 2024",
  "Tue Nov 12 17:42:58 EST 2024",
  "Tue Nov 12 17:43:07 EST 2024",
  "Tue Nov 12 17:43:16 EST 2024",
  "Tue Nov 12 17:43:26 EST 2024",
  "Tue Nov 12 17:43:36 EST 2024",
  "Tue Nov 12 17:43:48 EST 2024",
  "Tue Nov 12 17:43:58 EST 2024",
  "Tue Nov 12 17:44:07 EST 2024",
  "Tue Nov 12 17:44:16 EST 2024",
  "Tue Nov 12 17:44:25 EST 2024",
  "Tue Nov 12 17:44:35 EST 2024",
  "Tue Nov 12 17:44:47 EST 2024",
  "Tue Nov 12 17:44:56 EST 2024",
  "Tue Nov 12 17:45:05 EST 2024",
  "Tue Nov 12 17:45:15 EST 2024",
  "Tue Nov 12 17:45:24 EST 2024",
  "Tue Nov 12 17:45:33 EST 2024",
  "Tue Nov 12 17:45:46 EST 2024",
  "Tue Nov 12 17:45:55 EST 2024",
  "Tue Nov 12 17:46:05 EST 2024",
  "Tue Nov 12 17:46:14 EST 2024",
  "Tue Nov 12 17:46:24 EST 2024",
  "Tue Nov 12 17:46:33 EST 2024",
  "Tue Nov 12 17:46:45 EST 2024",
  "Tue Nov 12 17:46:54 EST 2024",
  "Tue Nov 12 17:47:04 EST 2024",
  "Tue Nov 12 17:47:13 EST 2024",
  "Tue Nov 12 17:47:23 EST 2024",
  "Tue Nov 12 17:47:32 EST 2024",
  "Tue Nov 12 17:47:44 EST 2024",
  "Tue Nov 12 17:47:53 EST 2024",
  "Tue Nov 12 17:48:03 EST 2024",
  "Tue Nov 12 17:48:13 EST 2024",
  "Tue Nov 12 17:48:22 EST 2024",
  "Tue Nov 12 17:48:32 EST 2024",
  "Tue Nov 12 17:48:44 EST 2024",
  "Tue Nov 12 17:48:54 EST 2024",
  "Tue Nov 12 17:49:04 EST 2024",
  "Tue Nov 12 17:49:15 EST 2024",
  "Tue Nov 12 17:49:25 EST 2024",
  "Tue Nov 12 17:49:35 EST 2024",
  "Tue Nov 12 17:49:50 EST 2024",
  "Tue Nov 12 17:49:59 EST 2024",
  "Tue Nov 12 17:50:08 EST 2024",
  "Tue Nov 12 17:50:17 EST 2024",
  "Tue Nov 12 17:50:28 EST 2024",
  "Tue Nov 12 17:50:38 EST 2024",
  "Tue Nov 12 17:50:47 EST 2024",
  "Tue Nov 12 17:50:57 EST 2024",
  "Tue Nov 12 17:51:07 EST 2024",
  "Tue Nov 12 17:51:16 EST 2024",
  "Tue Nov 12 17:51:26 EST 2024",
  "Tue Nov 12 17:51:35 EST 2024"
)

$rows = @()
for ($r = 2; $r -le 40; $r++) { $rows += $r }
for ($r = 44; $r -le 61; $r++) { $rows += $r }

for ($i = 0; $i -lt $rows.Length; $i++) {
  $ws.Cells.Item($rows[$i], 2).Value = $values[$i]
}
